$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.627.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.08%  '
$ws.Range("D3").Value = "'3.690.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +6.04%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'419.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").Value = "'130.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("D7").Value = "'3.678.84"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.98%  '
$ws.Range("D8").Value = "'0.644"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.86%  '
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").Value = "'0.759"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.19%  '
$ws.Range("D11").Value = "'0.181"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +10.43%  '
$ws.Range("D12").Value = "'0.0000392"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +47.62%  '
$ws.Range("D13").Value = "'42.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").Value = "'10.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.33%  '
$ws.Range("D15").Value = "'4.274.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.97%  '
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = "'20.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").Value = "'3.690.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.59%  '
$ws.Range("D19").Value = "'13.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.05%  '
$ws.Range("E20").Value = '  +2.63%  '
$ws.Range("D21").Value = "'66.632.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.25%  '
$ws.Range("D22").Value = "'443.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.04%  '
$ws.Range("D23").Value = "'16.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +21.17%  '
$ws.Range("D24").Value = "'89.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.22%  '
$ws.Range("D25").Value = "'3.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.32%  '
$ws.Range("D26").Value = "'37.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.04%  '
$ws.Range("D27").Value = "'10.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.85%  '
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("E29").Value = '  +4.50%  '
$ws.Range("D30").Value = "'12.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("D31").Value = "'0.123"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.75%  '
$ws.Range("D32").Value = "'2.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.73%  '
$ws.Range("E33").Value = '  -4.28%  '
$ws.Range("E34").Value = '  -0.92%  '
$ws.Range("D35").Value = "'41.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.45%  '
$ws.Range("D36").Value = "'57.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.40%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").Value = "'0.0494"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.19%  '
$ws.Range("D39").Value = "'3.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +37.83%  '
$ws.Range("D40").Value = "'0.0₃0744"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +13.25%  '
$ws.Range("E41").Value = '  +4.14%  '
$ws.Range("D42").Value = "'29.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +32.68%  '
$ws.Range("D43").Value = "'0.996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("E44").Value = '  +1.37%  '
$ws.Range("D45").Value = "'148.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.36%  '
$ws.Range("E46").Value = '  +3.65%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = "'2.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.20%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = "'4.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.95%  '
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").Value = "'2.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.70%  '
$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").Value = "'0.305"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.54%  '
$ws.Range("E51").Value = '  +15.22%  '